$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.975.52"
$ws.Range("D3").Value = "1.814.22"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'312.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.4289"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.3665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.07244"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +23.76%  "
$ws.Range("D11").Value = "'0.8620"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "'21.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'5.402"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "'6.599"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "'0.06938"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "'81.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'0.000008893"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'15.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "27.018.08"
$ws.Range("D22").Value = "'5.168"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "2.423.65"
$ws.Range("E23").Value = "  +21.95%  "
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'153.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "'1.879"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'18.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").Value = "'5.226"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("D29").Value = "'1.897"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.51%  "
$ws.Range("D30").Value = "'114.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'0.08935"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +6.87%  "
$ws.Range("D33").Value = "'0.7472"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").Value = "'4.414"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "'2.808"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("D38").Value = "'0.05202"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'0.01922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").Value = "'0.5096"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").Value = "'2.748"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").Value = "'6.481"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("D44").Value = "'8.336"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").Value = "'106.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").Value = "'10.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'0.4566"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "'1.642"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").Value = "'0.06208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'1.840"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.80%  "
